# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 159
$wsExhibit.Range("F4").Value = 573
$wsExhibit.Range("F5").Value = 1795
$wsExhibit.Range("F9").Value = 2183
$wsExhibit.Range("F10").Value = 47
$wsExhibit.Range("F12").Value = 1372
$wsExhibit.Range("F13").Value = 482
$wsExhibit.Range("F14").Value = 27
$wsExhibit.Range("F23").Value = 1185
$wsExhibit.Range("F25").Value = 352
$wsExhibit.Range("F26").Value = 179
$wsExhibit.Range("F27").Value = 277
$wsExhibit.Range("F28").Value = 346

# Sheet "全部类型" (sheet4): row -> new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 159
$wsAll.Range("F4").Value = 573
$wsAll.Range("F5").Value = 1795
$wsAll.Range("F10").Value = 2183
$wsAll.Range("F11").Value = 47
$wsAll.Range("F13").Value = 1372
$wsAll.Range("F14").Value = 482
$wsAll.Range("F15").Value = 27
$wsAll.Range("F24").Value = 1185
$wsAll.Range("F26").Value = 352
$wsAll.Range("F27").Value = 179
$wsAll.Range("F28").Value = 277
$wsAll.Range("F29").Value = 346
